# Insert a new data row at row 366 (pushing the existing rows 366-477 down to
# 367-478) on the active sheet, then populate the new row with its values.
#
# This corresponds to the diff, which shows every row from 366 through 477
# taking on the values previously held by the row immediately above it
# (i.e. a new record was inserted before the former row 366), and a brand
# new trailing row (478) appearing with the data that used to be in row 477.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 366..477 down to 367..478, leaving row 366 blank (but carrying
# the date-format style already used by the rest of column D).
$ws.Rows.Item(366).Insert()

# Populate the newly-inserted row 366 with the new record's values.
$ws.Cells.Item(366, 1).Value = 4
$ws.Cells.Item(366, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(366, 3).Value = "Los Lagos"
$ws.Cells.Item(366, 4).Value = 45215
$ws.Cells.Item(366, 5).Value = 10
$ws.Cells.Item(366, 6).Value = 100112032
$ws.Cells.Item(366, 7).Value = "Zapallo italiano"
$ws.Cells.Item(366, 8).Value = "Sin especificar"
$ws.Cells.Item(366, 9).Value = "Primera"
$ws.Cells.Item(366, 10).Value = 70
$ws.Cells.Item(366, 11).Value = 22000
$ws.Cells.Item(366, 12).Value = 22000
$ws.Cells.Item(366, 13).Value = 22000
$ws.Cells.Item(366, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(366, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(366, 16).Value = 440
$ws.Cells.Item(366, 17).Value = 50
$ws.Cells.Item(366, 18).Value = "Hortaliza"
